# Update scraped "想去人数" (interest count) figures for the regenerated
# gh-pages data snapshot (commit 456a3b4).
#
# The same events are listed on both the "展览" (Exhibition) sheet and the
# "全部类型" (All types) sheet, so the figures must be updated in both places.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> cell address -> new value
$updates = @{
    "展览" = @{
        "F2"  = 4278
        "F3"  = 2429
        "F6"  = 44
        "F7"  = 55
        "F9"  = 124
        "F10" = 133
        "F11" = 152
        "F12" = 1586
        "F13" = 292
        "F14" = 3306
    }
    "全部类型" = @{
        "F2"  = 4278
        "F3"  = 2429
        "F7"  = 44
        "F8"  = 55
        "F11" = 124
        "F12" = 133
        "F13" = 152
        "F16" = 1586
        "F17" = 292
        "F18" = 3306
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellAddr in $cellUpdates.Keys) {
        $ws.Range($cellAddr).Value = $cellUpdates[$cellAddr]
    }
}
